{"js": "const replacements = [\n  [\"520\u00f75=104, 0\", \"623\u00f78=77, 7\"],\n  [\"785\u00f73=261, 2\", \"756\u00f75=151, 1\"],\n  [\"197\u00f79=21, 8\", \"730\u00f78=91, 2\"],\n  [\"957\u00f73=319, 0\", \"397\u00f78=49, 5\"],\n  [\"956\u00f78=119, 4\", \"938\u00f73=312, 2\"],\n  [\"625\u00f75=125, 0\", \"377\u00f73=125, 2\"],\n  [\"760\u00f79=84, 4\", \"964\u00f78=120, 4\"],\n  [\"711\u00f78=88, 7\", \"516\u00f79=57, 3\"],\n  [\"139\u00f77=19, 6\", \"259\u00f78=32, 3\"],\n  [\"224\u00f77=32, 0\", \"711\u00f75=142, 1\"],\n  [\"669\u00f79=74, 3\", \"544\u00f75=108, 4\"],\n  [\"598\u00f74=149, 2\", \"747\u00f72=373, 1\"],\n  [\"792\u00f78=99, 0\", \"534\u00f78=66, 6\"],\n  [\"741\u00f73=247, 0\", \"895\u00f76=149, 1\"],\n  [\"228\u00f78=28, 4\", \"900\u00f73=300, 0\"],\n  [\"296\u00f73=98, 2\", \"920\u00f75=184, 0\"],\n  [\"450\u00f74=112, 2\", \"449\u00f75=89, 4\"],\n  [\"989\u00f76=164, 5\", \"854\u00f73=284, 2\"],\n  [\"516\u00f73=172, 0\", \"463\u00f79=51, 4\"],\n  [\"939\u00f73=313, 0\", \"625\u00f78=78, 1\"],\n  [\"462\u00f74=115, 2\", \"350\u00f77=50, 0\"],\n  [\"703\u00f74=175, 3\", \"156\u00f73=52, 0\"],\n  [\"118\u00f73=39, 1\", \"192\u00f72=96, 0\"],\n  [\"147\u00f78=18, 3\", \"685\u00f79=76, 1\"],\n  [\"106\u00f72=53, 0\", \"360\u00f79=40, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"520\u00f75=104, 0\", \"623\u00f78=77, 7\")\n  ,@(\"785\u00f73=261, 2\", \"756\u00f75=151, 1\")\n  ,@(\"197\u00f79=21, 8\", \"730\u00f78=91, 2\")\n  ,@(\"957\u00f73=319, 0\", \"397\u00f78=49, 5\")\n  ,@(\"956\u00f78=119, 4\", \"938\u00f73=312, 2\")\n  ,@(\"625\u00f75=125, 0\", \"377\u00f73=125, 2\")\n  ,@(\"760\u00f79=84, 4\", \"964\u00f78=120, 4\")\n  ,@(\"711\u00f78=88, 7\", \"516\u00f79=57, 3\")\n  ,@(\"139\u00f77=19, 6\", \"259\u00f78=32, 3\")\n  ,@(\"224\u00f77=32, 0\", \"711\u00f75=142, 1\")\n  ,@(\"669\u00f79=74, 3\", \"544\u00f75=108, 4\")\n  ,@(\"598\u00f74=149, 2\", \"747\u00f72=373, 1\")\n  ,@(\"792\u00f78=99, 0\", \"534\u00f78=66, 6\")\n  ,@(\"741\u00f73=247, 0\", \"895\u00f76=149, 1\")\n  ,@(\"228\u00f78=28, 4\", \"900\u00f73=300, 0\")\n  ,@(\"296\u00f73=98, 2\", \"920\u00f75=184, 0\")\n  ,@(\"450\u00f74=112, 2\", \"449\u00f75=89, 4\")\n  ,@(\"989\u00f76=164, 5\", \"854\u00f73=284, 2\")\n  ,@(\"516\u00f73=172, 0\", \"463\u00f79=51, 4\")\n  ,@(\"939\u00f73=313, 0\", \"625\u00f78=78, 1\")\n  ,@(\"462\u00f74=115, 2\", \"350\u00f77=50, 0\")\n  ,@(\"703\u00f74=175, 3\", \"156\u00f73=52, 0\")\n  ,@(\"118\u00f73=39, 1\", \"192\u00f72=96, 0\")\n  ,@(\"147\u00f78=18, 3\", \"685\u00f79=76, 1\")\n  ,@(\"106\u00f72=53, 0\", \"360\u00f79=40, 0\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute(\n    [ref]$oldText,\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]$newText,\n    [ref]2\n  )\n}"}
